# gx_lotmapping/std_curve_params.xlsx — 2024-04-26 lot mapping update
#
# 1) Left-align the "Cartridge Lot Number" column for the existing
#    numeric-lot rows (A48:A59) — matches the new cellXfs style that
#    appears in the target workbook (applyAlignment horizontal="left").
# 2) Append a new lot-mapping row (row 60) for cartridge lot 43302,
#    reusing the same curve IDs as the most recently added GX23-0002
#    SARS-CoV-2 lots (row 59).
# 3) Leave the sheet's selection parked where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Left-align cartridge lot numbers already present in A48:A59 ---
$ws.Range("A48:A59").HorizontalAlignment = -4131   # xlLeft

# --- 2) New row 60: cartridge lot 43302 ---
$ws.Range("A60").Value = 43302
$ws.Range("B60").Value = "GX23-0002_SARS_COV_2"
$ws.Range("C60").Value = "GX23-0001_FluA1"
$ws.Range("D60").Value = "GX23-0001_FluA2"
$ws.Range("E60").Value = "GX23-0001_FluB"
$ws.Range("F60").Value = "GX23-0001_RSV"

# --- 3) Restore the author's last on-screen selection ---
$ws.Range("D67").Select()
